$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 12775
$ws1.Range("F5").Value = 78
$ws1.Range("F6").Value = 56
$ws1.Range("F10").Value = 12678
$ws1.Range("F11").Value = 265
$ws1.Range("F13").Value = 7434
$ws1.Range("F14").Value = 7446
$ws1.Range("F15").Value = 178
$ws1.Range("F16").Value = 84
$ws1.Range("F18").Value = 117
$ws1.Range("F19").Value = 974
$ws1.Range("F23").Value = 181
$ws1.Range("F24").Value = 10

# Sheet "全部类型" (index 4) - column F ("想去人数") updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 12775
$ws4.Range("F6").Value = 78
$ws4.Range("F7").Value = 56
$ws4.Range("F11").Value = 12678
$ws4.Range("F12").Value = 265
$ws4.Range("F14").Value = 7438
$ws4.Range("F15").Value = 7448
$ws4.Range("F16").Value = 178
$ws4.Range("F17").Value = 84
$ws4.Range("F19").Value = 117
$ws4.Range("F20").Value = 974
$ws4.Range("F25").Value = 181
$ws4.Range("F26").Value = 10
